$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1849.5
$ws.Range("I4").Value = 699
$ws.Range("K4").Value = 699
$ws.Range("M4").Value = -585

$ws.Range("H19").Value = 1148.8
$ws.Range("I19").Value = 633
$ws.Range("J19").Value = 1369.8572
$ws.Range("K19").Value = 633
$ws.Range("L19").Value = 1369.8572
$ws.Range("M19").Value = -458
$ws.Range("N19").Value = -1719.8572

$ws.Range("H138").Value = 3147.25

$ws.Range("H141").Value = 2480.7693
$ws.Range("I141").Value = 2402
$ws.Range("J141").Value = 2658
$ws.Range("K141").Value = 7206
$ws.Range("L141").Value = 7974
$ws.Range("M141").Value = -2026
$ws.Range("N141").Value = -18334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3082689
$ws.Range("I32").Value = 3046384.8
$ws.Range("J32").Value = 3500187.5
$ws.Range("K32").Value = 3046384.8
$ws.Range("L32").Value = 3500187.5
$ws.Range("M32").Value = -3046097.8
$ws.Range("N32").Value = -3500761.5

$ws.Range("H36").Value = 13955.889
$ws.Range("I36").Value = 15679.167
$ws.Range("J36").Value = 10509.333
$ws.Range("K36").Value = 15679.167
$ws.Range("L36").Value = 10509.333
$ws.Range("M36").Value = -15333.167
$ws.Range("N36").Value = -11201.333

$ws.Range("H103").Value = 30500
$ws.Range("J103").Value = 30500
$ws.Range("L103").Value = 30500
$ws.Range("N103").Value = -32844

$ws.Range("H110").Value = 1436.125
$ws.Range("I110").Value = 1496.5714
$ws.Range("K110").Value = 1496.5714
$ws.Range("M110").Value = 548.4286

$ws.Range("H122").Value = 12808.182
$ws.Range("I122").Value = 16788.875
$ws.Range("J122").Value = 2193
$ws.Range("K122").Value = 50366.625
$ws.Range("L122").Value = 6579
$ws.Range("M122").Value = -47916.625
$ws.Range("N122").Value = -11479

$ws.Range("H132").Value = 1205.5714
$ws.Range("I132").Value = 1185.75
$ws.Range("K132").Value = 3557.25
$ws.Range("M132").Value = -1027.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H106").Value = 19557
$ws.Range("J106").Value = 19557
$ws.Range("L106").Value = 19557
$ws.Range("N106").Value = -22081

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 170000
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

$ws.Range("H104").Value = 49874.5
$ws.Range("J104").Value = 49874.5
$ws.Range("L104").Value = 49874.5
$ws.Range("N104").Value = -55116.5

$ws.Range("H122").Value = 673.3
$ws.Range("J122").Value = 789.5
$ws.Range("L122").Value = 2368.5
$ws.Range("N122").Value = -7268.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 725
$ws.Range("J17").Value = 1000
$ws.Range("L17").Value = 3000
$ws.Range("N17").Value = -3338

$ws.Range("H34").Value = 1336.5
$ws.Range("J34").Value = 1327.7142
$ws.Range("L34").Value = 3983.1426
$ws.Range("N34").Value = -4151.142599999999

$ws.Range("H39").Value = 732.8
$ws.Range("J39").Value = 1000
$ws.Range("L39").Value = 3000
$ws.Range("N39").Value = -3588

$ws.Range("H55").Value = 1000
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 3000
$ws.Range("M55").Value = -2823
$ws.Range("N55").Value = -3354

$ws.Range("H68").Value = 2676.3044
$ws.Range("I68").Value = 1115.4
$ws.Range("J68").Value = 2866.6584
$ws.Range("K68").Value = 3346.2
$ws.Range("L68").Value = 8599.975199999999
$ws.Range("M68").Value = -2535.2
$ws.Range("N68").Value = -10221.9752

$ws.Range("H71").Value = 2676.3044
$ws.Range("I71").Value = 1115.4
$ws.Range("J71").Value = 2866.6584
$ws.Range("K71").Value = 10038.6
$ws.Range("L71").Value = 25799.9256
$ws.Range("M71").Value = -5982.6
$ws.Range("N71").Value = -33911.9256

$ws.Range("H140").Value = 5414.923
$ws.Range("I140").Value = 1474.25
$ws.Range("K140").Value = 4422.75
$ws.Range("M140").Value = 757.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1678.0769
$ws.Range("I122").Value = 1501.4546
$ws.Range("K122").Value = 4504.3638
$ws.Range("M122").Value = -2054.3638

$ws.Range("H132").Value = 3584.8125
$ws.Range("I132").Value = 2194.2
$ws.Range("K132").Value = 6582.599999999999
$ws.Range("M132").Value = -4052.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 344.66666
$ws.Range("I16").Value = 305.77777
$ws.Range("J16").Value = 461.33334
$ws.Range("K16").Value = 305.77777
$ws.Range("L16").Value = 461.33334
$ws.Range("M16").Value = -135.77777
$ws.Range("N16").Value = -801.33334

$ws.Range("H22").Value = 1391.5
$ws.Range("I22").Value = 774.5
$ws.Range("J22").Value = 1700
$ws.Range("K22").Value = 774.5
$ws.Range("L22").Value = 1700
$ws.Range("M22").Value = -479.5
$ws.Range("N22").Value = -2290

$ws.Range("H27").Value = 1391.5
$ws.Range("I27").Value = 774.5
$ws.Range("J27").Value = 1700
$ws.Range("K27").Value = 774.5
$ws.Range("L27").Value = 1700
$ws.Range("M27").Value = -667.5
$ws.Range("N27").Value = -1914

$ws.Range("H61").Value = 1675.4445
$ws.Range("J61").Value = 1795.2
$ws.Range("L61").Value = 1795.2
$ws.Range("N61").Value = -2199.2

$ws.Range("H113").Value = 1675.4445
$ws.Range("J113").Value = 1795.2
$ws.Range("L113").Value = 1795.2
$ws.Range("N113").Value = -6135.2

$ws.Range("H132").Value = 2767.158
$ws.Range("I132").Value = 2831.8667
$ws.Range("K132").Value = 8495.6001
$ws.Range("M132").Value = -5965.6001

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 9000
$ws.Range("J5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("N5").Value = -9224

$ws.Range("H74").Value = 20590.334
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 20590.334
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 20590.334
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -22462.334

$ws.Range("H77").Value = 20590.334
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 20590.334
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 61771.00199999999
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -71131.00199999999

$ws.Range("H104").Value = 21993
$ws.Range("J104").Value = 21993
$ws.Range("L104").Value = 21993
$ws.Range("N104").Value = -28981

$ws.Range("H122").Value = 4236.778
$ws.Range("I122").Value = 2142.5
$ws.Range("J122").Value = 5912.2
$ws.Range("K122").Value = 6427.5
$ws.Range("L122").Value = 17736.6
$ws.Range("M122").Value = -3977.5
$ws.Range("N122").Value = -22636.6

$ws.Range("H126").Value = 3955.5557
$ws.Range("J126").Value = 5332.8335
$ws.Range("L126").Value = 15998.5005
$ws.Range("N126").Value = -20938.5005

$ws.Range("H136").Value = 2060.25
$ws.Range("I136").Value = 2069.7036
$ws.Range("K136").Value = 6209.110799999999
$ws.Range("M136").Value = -3659.110799999999
